$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter hours for "Work (misc.)" row (row 13): 2 hours on Thursday (F), 4 hours on Saturday (H)
$ws.Range("F13").Value = 2
$ws.Range("H13").Value = 4

# Move selection to H13 to match where the user left off
$ws.Range("H13").Select()
